$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new header row is inserted above the existing data (rows 1-20 -> 2-21).
$ws.Rows.Item(1).Insert()

# Fill the new header row. Values are written in this particular order
# (A, C, B, D) so that the shared-string table ends up with the same
# insertion order as the target workbook: item, escalas, pregunta,
# posibles respuestas.
$ws.Range("A1").Value = "item"
$ws.Range("C1").Value = "escalas"
$ws.Range("B1").Value = "pregunta"
$ws.Range("D1").Value = "posibles respuestas"

# Restore the view: scrolled so column C is leftmost, with D2 selected.
try {
    $excel.ActiveWindow.ScrollColumn = 3
} catch {
}
$ws.Range("D2").Select()
